$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.619.45"
$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").Value = "1.855.64"
$ws.Range("E3").Value = "  +0.38%  "

$ws.Range("D5").Value = "'264.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.93%  "

$ws.Range("D6").Value = "'0.9997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("E7").Value = "  -0.58%  "

$ws.Range("D8").Value = "'0.3285"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.36%  "

$ws.Range("D9").Value = "'0.06807"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.81%  "

$ws.Range("D10").Value = "'18.86"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.04%  "

$ws.Range("D11").Value = "'0.7775"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.05%  "

$ws.Range("D12").Value = "'0.07770"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.58%  "

$ws.Range("D13").Value = "1.840.45"
$ws.Range("E13").Value = "  -0.38%  "

$ws.Range("D14").Value = "'88.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.01%  "

$ws.Range("D15").Value = "'5.030"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.54%  "

$ws.Range("D16").Value = "'0.9995"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.07%  "

$ws.Range("D17").Value = "'14.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.07%  "

$ws.Range("D18").Value = "'0.000007977"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.83%  "

$ws.Range("D19").Value = "'0.9999"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("D20").Value = "26.632.87"
$ws.Range("E20").Value = "  +0.19%  "

$ws.Range("D21").Value = "2.084.14"
$ws.Range("E21").Value = "  +0.48%  "

$ws.Range("D22").Value = "'4.643"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.62%  "

$ws.Range("D23").Value = "'9.571"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.70%  "

$ws.Range("D24").Value = "'6.005"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.09%  "

$ws.Range("D25").Value = "'144.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.05%  "

$ws.Range("E26").Value = "  -6.86%  "

$ws.Range("D27").Value = "'1.672"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.46%  "

$ws.Range("D28").Value = "'17.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.10%  "

$ws.Range("D29").Value = "'112.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.90%  "

$ws.Range("D30").Value = "'4.208"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.68%  "

$ws.Range("D31").Value = "'4.160"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.87%  "

$ws.Range("D32").Value = "'0.08763"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.31%  "

$ws.Range("D33").Value = "'0.04840"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.39%  "

$ws.Range("D34").Value = "'1.141"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.42%  "

$ws.Range("D35").Value = "'0.7172"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.99%  "

$ws.Range("E36").Value = "  -0.07%  "

$ws.Range("D37").Value = "'3.117"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.25%  "

$ws.Range("D38").Value = "'0.01789"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.41%  "

$ws.Range("D39").Value = "'2.221"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.07%  "

$ws.Range("D40").Value = "'0.4896"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.49%  "

$ws.Range("D41").Value = "'113.03"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.92%  "

$ws.Range("D42").Value = "'0.9044"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.03%  "

$ws.Range("D43").Value = "'6.098"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.35%  "

$ws.Range("D44").Value = "'0.9997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("D45").Value = "'7.747"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.03%  "

$ws.Range("D46").Value = "'0.4203"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.34%  "

$ws.Range("D47").Value = "'9.153"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.29%  "

$ws.Range("D48").Value = "'0.05941"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.25%  "

$ws.Range("E49").Value = "  -3.92%  "

$ws.Range("D50").Value = "'35.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.05%  "

$ws.Range("D51").Value = "'0.8869"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.35%  "
